$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reuse the formatting (fill/alignment/wrap) of the row above (row 9) ---
$ws.Range("A9:D9").Copy() | Out-Null
$ws.Range("A10:D10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row height matches the other data rows (225) ---
$ws.Rows.Item(10).RowHeight = 225

# --- A10: new test file name ---
$ws.Range("A10").Value = "tb_amp_squid_top_test_variant_debug04.json"

# --- B10: parameters summary, rich text (red "8" sample-size run) ---
$ws.Range("B10").Value = "`"nb_sample_by_pixel`": 8,`r`n`"nb_pixel_by_frame`": 1,`r`n`"nb_frame_by_pulse`": 2048,`r`n`"nb_pulse`": 8,`r`n`"amp_squid_offset_correction`"/ « mode »:0,`r`n`"amp_squid_offset_correction`"/ «min_value »:0,`r`n`"amp_squid_offset_correction`"/ « max_value »:0,`r`n« pixel_result »/ « mode »:0,`r`n« pixel_result »/ « min_value »:0,`r`n« pixel_result »/ « max_value »:131071"
$b10r2 = $ws.Range("B10").Characters(23, 1)
$b10r2.Font.Size = 11
$b10r2.Font.Name = "Calibri"
$b10r2.Font.Color = 255
$b10r3 = $ws.Range("B10").Characters(24, 318)
$b10r3.Font.Size = 11
$b10r3.Font.Name = "Calibri"
$b10r3.Font.Color = 0

# --- C10: reading speed comment, rich text (red header line) ---
$ws.Range("C10").Value = "Continuous data valid`r`nContinuous ram1 check`r`n"
$c10r1 = $ws.Range("C10").Characters(1, 23)
$c10r1.Font.Size = 10
$c10r1.Font.Name = "Liberation Sans"
$c10r1.Font.Color = 255
$c10r2 = $ws.Range("C10").Characters(24, 23)
$c10r2.Font.Size = 11
$c10r2.Font.Name = "Calibri"
$c10r2.Font.Color = 0

# --- D10: comment ---
$ws.Range("D10").Value = "Auto-check`r`ndisable adc data (=0) to read the memory content`r`n1 pixel by frame"

# --- Selection / view state matches where the author left off editing ---
$ws.Range("D10").Select()
